# Update workbook to reflect carjacking data as of 2022-07-13
# (previously "through July 12").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab: "Through 2022-07-12" -> "Through 2022-07-13"
$ws.Name = "Through 2022-07-13"

# Update the column header label to match (shared string used only by B1)
$ws.Range("B1").Value2 = "July 2022 (through July 13)"

# Apply the per-neighborhood / per-month count updates.
# Each of these is a "July" column for some year, bumped because a
# carjacking that occurred on 2022-07-13 is now included in the tally
# (for current year it directly increments; the historical "July" columns
# for prior years reflect the same additional day now being counted).
$ws.Range("I2").Value2  = 7
$ws.Range("W2").Value2  = 2

$ws.Range("B3").Value2  = 6
$ws.Range("I3").Value2  = 2

$ws.Range("B5").Value2  = 3
$ws.Range("P5").Value2  = 4
$ws.Range("W5").Value2  = 2
$ws.Range("AK5").Value2 = 1
$ws.Range("AR5").Value2 = 4
$ws.Range("AY5").Value2 = 2

$ws.Range("AR8").Value2 = 3

$ws.Range("I14").Value2 = 1

$ws.Range("I19").Value2 = 3
$ws.Range("W19").Value2 = 2

$ws.Range("W27").Value2 = 1

$ws.Range("I29").Value2 = 3
$ws.Range("P29").Value2 = 2

$ws.Range("W30").Value2 = 1

$ws.Range("W35").Value2 = 1

$ws.Range("W47").Value2 = 1

$ws.Range("AR58").Value2 = 2

$ws.Range("I65").Value2 = 2

$ws.Range("B70").Value2 = 2

$ws.Range("AD89").Value2 = 1
